$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($ws, $addr, $val)
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 2
Set-TextCell $ws 'D2' '59.843.14'
$ws.Range('E2').Value = '  +1.01%  '

# Row 3
Set-TextCell $ws 'D3' '2.568.96'
$ws.Range('E3').Value = '  +1.77%  '

# Row 4
$ws.Range('E4').Value = '  -0.21%  '

# Row 5
Set-TextCell $ws 'D5' '503.63'
$ws.Range('E5').Value = '  -0.78%  '

# Row 6
Set-TextCell $ws 'D6' '152.03'
$ws.Range('E6').Value = '  -6.12%  '

# Row 7
$ws.Range('E7').Value = '  +0.65%  '

# Row 8
Set-TextCell $ws 'D8' '0.577'
$ws.Range('E8').Value = '  -5.53%  '

# Row 9
Set-TextCell $ws 'D9' '2.570.90'
$ws.Range('E9').Value = '  +0.53%  '

# Row 10
Set-TextCell $ws 'D10' '6.77'
$ws.Range('E10').Value = '  +7.45%  '

# Row 11
$ws.Range('E11').Value = '  -1.96%  '

# Row 12
Set-TextCell $ws 'D12' '0.343'
$ws.Range('E12').Value = '  +0.00%  '

# Row 13
$ws.Range('E13').Value = '  +0.34%  '

# Row 14
Set-TextCell $ws 'D14' '3.015.16'
$ws.Range('E14').Value = '  +2.16%  '

# Row 15
Set-TextCell $ws 'D15' '59.859.62'
$ws.Range('E15').Value = '  +1.29%  '

# Row 16
Set-TextCell $ws 'D16' '21.50'
$ws.Range('E16').Value = '  -3.05%  '

# Row 17
Set-TextCell $ws 'D17' '0.0000139'
$ws.Range('E17').Value = '  -0.36%  '

# Row 18
Set-TextCell $ws 'D18' '2.567.87'
$ws.Range('E18').Value = '  +0.56%  '

# Row 19
Set-TextCell $ws 'D19' '4.78'
$ws.Range('E19').Value = '  +0.43%  '

# Row 20
Set-TextCell $ws 'D20' '344.99'
$ws.Range('E20').Value = '  +2.92%  '

# Row 21
Set-TextCell $ws 'D21' '10.24'
$ws.Range('E21').Value = '  -0.46%  '

# Row 22
Set-TextCell $ws 'D22' '6.02'
$ws.Range('E22').Value = '  -0.65%  '

# Row 23
Set-TextCell $ws 'D23' '0.998'
$ws.Range('E23').Value = '  -0.25%  '

# Row 24
$ws.Range('E24').Value = '  +0.51%  '

# Row 25
Set-TextCell $ws 'D25' '0.417'
$ws.Range('E25').Value = '  +0.03%  '

# Row 26
$ws.Range('E26').Value = '  -2.87%  '

# Row 27
Set-TextCell $ws 'D27' '2.670.41'
$ws.Range('E27').Value = '  +3.29%  '

# Row 28
$ws.Range('E28').Value = '  +0.53%  '

# Row 29
Set-TextCell $ws 'D29' '0.0₃0841'
$ws.Range('E29').Value = '  +2.28%  '

# Row 30
Set-TextCell $ws 'D30' '7.43'
$ws.Range('E30').Value = '  -1.95%  '

# Row 31
$ws.Range('E31').Value = '  +0.23%  '

# Row 32
Set-TextCell $ws 'D32' '154.83'
$ws.Range('E32').Value = '  +0.09%  '

# Row 33
Set-TextCell $ws 'D33' '19.19'
$ws.Range('E33').Value = '  -1.24%  '

# Row 34
Set-TextCell $ws 'D34' '1.55'
$ws.Range('E34').Value = '  -1.22%  '

# Row 35
Set-TextCell $ws 'D35' '5.70'
$ws.Range('E35').Value = '  +3.28%  '

# Row 36
Set-TextCell $ws 'D36' '3.97'
$ws.Range('E36').Value = '  +0.97%  '

# Row 37
$ws.Range('E37').Value = '  -0.76%  '

# Row 38
Set-TextCell $ws 'D38' '0.852'
$ws.Range('E38').Value = '  +19.80%  '

# Row 39
$ws.Range('B39').Value = 'Fetch.AI'
$ws.Range('C39').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextCell $ws 'D39' '0.840'
$ws.Range('E39').Value = '  -2.85%  '

# Row 40
$ws.Range('B40').Value = 'Filecoin'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextCell $ws 'D40' '3.75'
$ws.Range('E40').Value = '  -0.11%  '

# Row 41
$ws.Range('E41').Value = '  +0.21%  '

# Row 42
$ws.Range('B42').Value = 'Bittensor'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextCell $ws 'D42' '297.01'
$ws.Range('E42').Value = '  +1.83%  '

# Row 43
$ws.Range('B43').Value = 'OKB'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextCell $ws 'D43' '35.38'
$ws.Range('E43').Value = '  +1.55%  '

# Row 44
$ws.Range('B44').Value = 'Stellar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextCell $ws 'D44' '0.0994'
$ws.Range('E44').Value = '  -1.72%  '

# Row 45
$ws.Range('B45').Value = 'Hedera'
$ws.Range('C45').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextCell $ws 'D45' '0.0558'
$ws.Range('E45').Value = '  +0.08%  '

# Row 46
$ws.Range('B46').Value = 'FirstDigitalUSD'
$ws.Range('C46').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextCell $ws 'D46' '0.998'
$ws.Range('E46').Value = '  +0.84%  '

# Row 47
Set-TextCell $ws 'D47' '0.612'
$ws.Range('E47').Value = '  -1.46%  '

# Row 48
Set-TextCell $ws 'D48' '19.61'
$ws.Range('E48').Value = '  +4.43%  '

# Row 49
Set-TextCell $ws 'D49' '4.84'
$ws.Range('E49').Value = '  -0.28%  '

# Row 50
Set-TextCell $ws 'D50' '0.0233'
$ws.Range('E50').Value = '  -2.25%  '

# Row 51
$ws.Range('B51').Value = 'Maker'
$ws.Range('C51').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextCell $ws 'D51' '2.006.40'
$ws.Range('E51').Value = '  +2.53%  '
